$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.658.28'
$ws.Range('D3').Value = '2.999.80'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.38'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.40%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.64%  '
$ws.Range('D9').Value = '2.998.97'
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('E10').Value = '  -2.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.89'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.463'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000228'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.35'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.02%  '
$ws.Range('E15').Value = '  +2.63%  '
$ws.Range('D16').Value = '3.495.25'
$ws.Range('E16').Value = '  -1.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.01'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').Value = '61.602.67'
$ws.Range('E18').Value = '  -1.81%  '
$ws.Range('D19').Value = '2.998.37'
$ws.Range('E19').Value = '  -1.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '455.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.03'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.689'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.59%  '
$ws.Range('E23').Value = '  -0.85%  '
$ws.Range('E24').Value = '  +1.67%  '
$ws.Range('E25').Value = '  -7.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.17'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.47'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('E29').Value = '  +1.75%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.01'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.52%  '
$ws.Range('E32').Value = '  -4.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.37'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('E34').Value = '  -0.38%  '
$ws.Range('D35').Value = '0.0₃0818'
$ws.Range('E35').Value = '  +2.71%  '
$ws.Range('E36').Value = '  -1.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.76'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('E38').Value = '  -3.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '9.19'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.31'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.123'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.12%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.89'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '400.30'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '39.30'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0353'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.80%  '
$ws.Range('E46').Value = '  -5.40%  '
$ws.Range('D47').Value = '2.719.79'
$ws.Range('E47').Value = '  -3.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.26'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.34%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.107'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.17'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.23%  '
